$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Party Master in Transation IN form"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "In progress"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Party Master in Transation Out form"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "In progress"

$ws.Range("F10").Select()
